$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in new burndown values for this update
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 4

$ws.Range("C12").Value = "Eman"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 6

# Recalculate totals (formulas in row 42 will pick these up automatically)
$excel.CalculateFullRebuild()

# Move the active selection to reflect where the user left off editing
$ws.Range("F13").Select() | Out-Null
